$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 18, pushing existing rows 18.. down by one.
$ws.Rows.Item(18).Insert()

# Populate the newly inserted row 18 with the new parameter definition.
$ws.Range("A18").Value = "general.maxNumberCompThreads"
$ws.Range("D18").Value = "structural_pipeline"
$ws.Range("E18").Value = "numeric"
$ws.Range("F18").Value = "scalar nonempty nonnegative"
$ws.Range("G18").Value = "standard"
$ws.Range("H18").Value = "Maximum number of computational threads used in pipeline. Value 0 lets MATLAB determine the most desirable number of computational threads (equal to the number of physical cores on the machine)."

# Match the selection left behind by the authoring session.
$ws.Range("A21").Select()
